$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '245.68'
$ws.Range('D2').Style = "Normal"

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '25.48'
$ws.Range('D3').Style = "Normal"

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.092'
$ws.Range('D4').Style = "Normal"

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.05570'
$ws.Range('D5').Style = "Normal"

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '6.478'
$ws.Range('D6').Style = "Normal"

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.017'
$ws.Range('D7').Style = "Normal"

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.8185'
$ws.Range('D8').Style = "Normal"

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.8437'
$ws.Range('D9').Style = "Normal"

# Row 10
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.009785'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '9OneONEBestin24h'

# Row 11
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.1342'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '10WazirXWRX'

# Row 12
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.06947'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.02859'
$ws.Range('D13').Style = "Normal"

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.09374'
$ws.Range('D14').Style = "Normal"

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.001519'
$ws.Range('D15').Style = "Normal"

# Row 16
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.006097'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '15TigerCashTCH'

# Row 17
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.499'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '16LEOLEO'

# Row 18
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.091'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '17BTSETokenBTSE'

# Row 19
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.3179'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '18BitpandaEcosystemTokenBEST'

# Row 20
$ws.Range('B20').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C20').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.03188'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '19LiechtensteinCryptoassetsExchangeLCX'

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.752'
$ws.Range('D22').Style = "Normal"

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.04710'
$ws.Range('D23').Style = "Normal"

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.001246'
$ws.Range('D25').Style = "Normal"

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.004632'
$ws.Range('D26').Style = "Normal"

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.00009701'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '26NitroExNTX'

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.03661'
$ws.Range('D40').Style = "Normal"

# Row 41
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.006211'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '40KickTokenKICK'

# Row 42
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1051'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '41BKEXTokenBKK'

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.002500'
$ws.Range('D43').Style = "Normal"

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.008303'
$ws.Range('D44').Style = "Normal"

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005304'
$ws.Range('D45').Style = "Normal"

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.1330'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '46CoinbaseStockTokenCOINWorstin24h'

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.002122'
$ws.Range('D48').Style = "Normal"

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00002100'
$ws.Range('D49').Style = "Normal"

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0002000'
$ws.Range('D50').Style = "Normal"

Write-Host "Applied all updates"